$wb = $excel.ActiveWorkbook

# Rename sheet "Phase 5" to "Phase5"
$ws = $wb.Worksheets.Item("Phase 5")
$ws.Name = "Phase5"

# Make sure this sheet is active/selected
$ws.Activate()
$ws.Select()

# Update B42 value from 0.02 to 0
$ws.Range("B42").Value = 0

# Update the view: scroll so A11 is the top-left cell, and select B43
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("B43").Select()
